# "Add files via upload" — refresh of the COVID-19 Valais daily figures table
# (sheet "Feuil1"). Updates the isolation/quarantine columns (N/O/P) for the
# previously-entered rows 181-204, corrects the new-positive-cases counts
# (C203/C204, which ripple into the cumulative B column through the existing
# shared formula), fills in the previously-blank row 205 (date 2020-09-17),
# and moves the frozen-pane scroll position / active selection down a week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column N (Nombre de cas en cours d'isolement) corrections, rows 181-199 ---
$ws.Range("N181").Value = 68
$ws.Range("N182").Value = 83
$ws.Range("N183").Value = 86
$ws.Range("N184").Value = 86
$ws.Range("N185").Value = 90
$ws.Range("N186").Value = 85
$ws.Range("N187").Value = 83
$ws.Range("N188").Value = 75
$ws.Range("N189").Value = 79
$ws.Range("N190").Value = 78
$ws.Range("N191").Value = 76
$ws.Range("N192").Value = 76
$ws.Range("N193").Value = 71
$ws.Range("N194").Value = 70
$ws.Range("N195").Value = 83
$ws.Range("N197").Value = 90
$ws.Range("N198").Value = 95
$ws.Range("N199").Value = 93

# --- Rows 200-204: N/O/P corrections, plus corrected new-case counts (C) ---
$ws.Range("N200").Value = 87
$ws.Range("O200").Value = 224

$ws.Range("N201").Value = 82
$ws.Range("O201").Value = 249
$ws.Range("P201").Value = 285

$ws.Range("N202").Value = 83
$ws.Range("O202").Value = 289
$ws.Range("P202").Value = 295

$ws.Range("C203").Value = 10
$ws.Range("N203").Value = 76
$ws.Range("O203").Value = 288
$ws.Range("P203").Value = 297

$ws.Range("C204").Value = 15
$ws.Range("N204").Value = 81
$ws.Range("O204").Value = 328
$ws.Range("P204").Value = 285

# --- Row 205 (2020-09-17) was entirely blank (formulas resolved to "");
#     fill in the day's figures so the shared formulas in B/H/J/K pick up a
#     real value. ---
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 8
$ws.Range("I205").Value = 1
$ws.Range("L205").Value = "0"
$ws.Range("M205").Value = "0"
$ws.Range("N205").Value = 71
$ws.Range("O205").Value = 287
$ws.Range("P205").Value = 259

# --- View state: scroll the frozen pane down one row and move the active
#     selection to B207 (bottom-right pane), matching where the sheet was
#     left after the new row was entered. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 195
$win.ScrollColumn = 2
$ws.Range("B207").Select()
